# Gangwon Juso data update
# - Row 2 (강원특별자치도영월의료원): corrected latitude/longitude
# - Row 15 (홍성국신경정신과의원): corrected address (street number) and longitude
# - Move active selection to F16
# - Resize the workbook window

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the latitude / longitude for row 2 (강원특별자치도영월의료원)
$ws.Range("E2").Value = 37.186540000000001
$ws.Range("F2").Value = 128.46510000000001

# Fix the street address and longitude for row 15 (홍성국신경정신과의원)
$ws.Range("C15").Value = "원주시 원일로 150-137.35219"
$ws.Range("F15").Value = 127.9472

# Move the selection to F16, matching the saved view state
[void]$ws.Range("F16").Select()

# Match the resized workbook window from the author's session
$win = $wb.Windows.Item(1)
$win.Width = 11420
